$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Ajo" (Chino, Primera) at
# Feria Lagunitas de Puerto Montt. Insert it as the new first data row of
# this block (row 96), pushing the existing rows 96-145 down to 97-146
# (the last of which, row 146, is simply the former row 145 shifted down).
$ws.Rows.Item(96).Insert()

$ws.Range("A96").Value = 4
$ws.Range("B96").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C96").Value = "Los Lagos"
$ws.Range("D96").Value = 44452
$ws.Range("E96").Value = 10
$ws.Range("F96").Value = 100112003
$ws.Range("G96").Value = "Ajo"
$ws.Range("H96").Value = "Chino"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 80
$ws.Range("K96").Value = 18000
$ws.Range("L96").Value = 18000
$ws.Range("M96").Value = 18000
$ws.Range("N96").Value = "$/caja 10 kilos"
$ws.Range("O96").Value = "China"
$ws.Range("P96").Value = 1800
$ws.Range("Q96").Value = 10
$ws.Range("R96").Value = "Hortaliza"
